$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 34: change the "Post Purchase Upsell Pricing" values from merged text
# strings into real numeric currency values (B34:M34 = 49.95, N34:Y34 = 29.95),
# and shrink the row height now that it's single line instead of wrapped text.
$ws.Range("B34:M34").Value = 49.95
$ws.Range("N34:Y34").Value = 29.95
$ws.Rows.Item(34).RowHeight = 15.75

# Row 44 ("Post Purchase Upsell Promotion 4"): swap the free-gift copy.
$ws.Range("B44:M44").Value = "Free Bold & Beautiful Mascara - WY1A0027"

# Update the view state left over from scrolling/selecting in the sheet.
$ws.Range("O34:Y34").Select()
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 19
